# Applies the personnel-change edits described in the commit
# "Solventar probelmas dce planeacion" to the BAJA document.

$d = $word.ActiveDocument

function Replace-Exact {
    param(
        [string]$OldText,
        [string]$NewText
    )
    $range = $d.Content
    $range.Find.Execute($OldText, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $NewText, 2)
}

# CATEGORÍA ACTUAL code
Replace-Exact "2A0508A" "2S0101A"

# Job title text
Replace-Exact "OFICIAL ADMINISTRATIVO 5A" "AUXILIAR 1A"

# DEPENDENCIA (department) text
Replace-Exact "CONTABILIDAD GUBERNAMENTAL" "CONTROL DE REC. HUMANOS Y SUELDOS APLICADOS"

# CLAVE DEL PROYECTO (project key)
Replace-Exact "1140031490300000120" "1140020000000000220"

# Apellido paterno - add accent
Replace-Exact "JIMENEZ" "JIMÉNEZ"

# Nombre(s) - add accent
Replace-Exact "ISIDRO NOE" "ISIDRO NOÉ"

# "A PARTIR DEL" day-of-month field (30 -> 03). MatchWholeWord (the 2nd
# argument to Replace-Exact's Find.Execute) keeps this from touching the
# "30" that appears as a substring of other numbers (e.g. the project key).
Replace-Exact "30" "03"
